$d = $word.ActiveDocument

function Insert-AtStart($paraIndex, $text) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $rZero = $d.Range($r.Start, $r.Start)
    $rZero.Text = $text
}

function Insert-AfterFirstChar($paraIndex, $text) {
    # Inserts $text immediately after the first character of the paragraph
    # (used to add text right after a lone <w:tab/> run) while keeping the
    # tab run itself intact.
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $rFirst = $d.Range($r.Start, $r.Start + 1)
    $rFirst.InsertAfter($text)
}

# --- Education block ---
# Tab run gains the date range text right after the tab character
# (done first, while the tab is still the paragraph's first character)
Insert-AfterFirstChar 4 "2023-01-20 - 2023-01-31"
# Bold (empty) run gains "k"
Insert-AtStart 4 "k"
# The next (empty) paragraph gains a new run with "kk"
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertAfter("kk")

# "Current GPA: " -> "Current GPA: 4"
$d.Content.Find.Execute("Current GPA: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Current GPA: 4", 2)

# Education bullet sentence rewrite
$d.Content.Find.Execute(
    "Revise the sentence to sound more powerful: Convince employers of your suitability for the job with your resume.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Achieved 880k in sales within one quarter, exceeding target by 8%.", 2)

# --- Experience block ---
# Tab run (paragraph 9) gains " - " right after the tab character
Insert-AfterFirstChar 9 " - "

# Experience bullet sentence rewrite
$d.Content.Find.Execute(
    "Revise the sentence to read:Craft a more compelling summary statement for your resume.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Revise the sentence to: Craft a persuasive argument for inclusion on your resume.", 2)

# --- CCA block ---
# Tab run (paragraph 13) gains " - " right after the tab character
Insert-AfterFirstChar 13 " - "

# CCA bullet sentence rewrite
$d.Content.Find.Execute(
    "Revise the sentence to: Strengthen the impact of this resume.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Revise the sentence to read: Strengthen your resume with a compelling narrative.", 2)

# --- Volunteer block ---
# Tab run (paragraph 17) gains " - " right after the tab character
Insert-AfterFirstChar 17 " - "

# Volunteer bullet sentence rewrite
$d.Content.Find.Execute(
    "Rewrite the sentence as:  Craft a persuasive resume that stands out.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Revise the sentence to read: Create a compelling resume to make a strong impression.", 2)
